# "update statistics by excel"
# The user extended the little 10-column "fund statistics" block (columns A:J,
# pattern copied from the last existing data row) by three more rows: two rows
# appended right after the existing block (rows 29 and 30), then a third one a
# couple of rows further down (row 33), leaving rows 31-32 empty - exactly the
# way you'd get by selecting the last row, copying it, and pasting it into the
# new rows with the fill handle / paste.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteAll = -4104

# The last existing row of the repeating A:J block.
$source = $ws.Range("A28:J28")
$source.Copy() | Out-Null

$ws.Range("A29:J29").PasteSpecial($xlPasteAll) | Out-Null
$ws.Range("A30:J30").PasteSpecial($xlPasteAll) | Out-Null
$ws.Range("A33:J33").PasteSpecial($xlPasteAll) | Out-Null

$excel.CutCopyMode = 0

# Leave the view scrolled down a bit and the newly added last row selected,
# like it was right after typing/pasting it in.
$ws.Range("A33:J33").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 13
